$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, shifting rows 12-13 down to 13-14
$ws.Rows("12:12").Insert()

# Fill in the new row 12 with the new weekly price entry
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 45001
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112032
$ws.Range("G12").Value = "Zapallo italiano"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 40
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = "$/caja 60 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 167
$ws.Range("Q12").Value = 60
$ws.Range("R12").Value = "Hortaliza"
